# CET350 - update tasklist
# 1) "Assignment:" number: 1 -> 2
# 2) Each "Date:" entry in the roster table: 01/25/22 -> 02/07/22
#    (typed as a selection-replace, which Word records as three runs:
#     "0" | "2/07" | "/22")

$d = $word.ActiveDocument

# --- 1) Assignment number -----------------------------------------------
$p1 = $d.Paragraphs(1).Range
$found = $p1.Find.Execute("1", $false, $true, $false, $false, $false, `
                           $true, 1, $false, "", 0)
if ($found) {
    $p1.Text = "2"
}

# --- 2) Roster dates ------------------------------------------------------
$searchFrom = 0
$docEnd = $d.Content.End

while ($true) {
    $rng = $d.Range($searchFrom, $docEnd)
    $found = $rng.Find.Execute("01/25/22", $false, $false, $false, $false, `
                                $false, $true, 1, $false, "", 0)
    if (-not $found) {
        break
    }

    $cellStart = $rng.Start

    # Re-typing the middle of the date (the selection a user would have
    # replaced) so Word keeps the untouched "0" prefix and "/22" suffix in
    # their own runs, with the newly typed "2/07" as a run in between.
    $prefix = $d.Range($cellStart, $cellStart + 1)
    $middle = $d.Range($cellStart + 1, $cellStart + 5)
    $suffix = $d.Range($cellStart + 5, $cellStart + 8)

    $middle.Text = "2/07"

    $middle2 = $d.Range($cellStart + 1, $cellStart + 5)
    $middle2.Bold = $true
    $middle2.Bold = $false

    $searchFrom = $cellStart + 8
    $docEnd = $d.Content.End
}
